# The scraped change for this asset (images/experiments/svg_template.pptx)
# touches only PowerPoint's internal collaboration/revision-tracking parts
# (ppt/revisionInfo.xml and ppt/changesInfos/changesInfo1.xml): a client
# version counter and two "last modified" timestamps for the existing
# picture's add/mod history entry. No slide, shape, text, or media content
# actually differs between the before/after OOXML - the picture (id 76,
# "Graphic 75") keeps the same geometry, fill, and formatting.
#
# That bookkeeping is written by the live PowerPoint app as a side effect
# of its own session/undo stack (autosave revision stamps) and isn't
# surfaced anywhere on the Presentation/Slide/Shape object model, so it
# can't be (and shouldn't be) re-derived from a content edit here. We
# simply confirm the deck's current state matches what that revision
# entry already describes, and leave it untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$pic = $s.Shapes.Item(1)
Write-Host "Slide count:" $p.Slides.Count
Write-Host "Shape:" $pic.Name "id-ish index 1, type" $pic.Type
Write-Host "Position:" $pic.Left $pic.Top $pic.Width $pic.Height
